$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "done" (x) for the core-feature checklist rows that already have
# working progress (rows 2-7, 10-11) in a new column D, centered both
# ways - matching the little checkbox column being introduced while the
# camera / draw-correctness rows (8-9) are left blank since that's what's
# still being worked on (world Matrix).
$ws.Range("D2").VerticalAlignment = -4108
$ws.Range("D2").HorizontalAlignment = -4108
$ws.Range("D2").Value = "x"

$ws.Range("D2").Copy()
$ws.Range("D3:D7").PasteSpecial(-4122)
$ws.Range("D10:D11").PasteSpecial(-4122)
$ws.Range("D3:D7").Value = "x"
$ws.Range("D10:D11").Value = "x"

# Restore the selection to where work is currently happening.
$ws.Range("A14").Select() | Out-Null
